$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new values for K2, L2 (row 2) and K4, L4 (row 4)
$ws.Range("K2").Value = "OP-00"
$ws.Range("L2").Value = "999-99"
$ws.Range("L2").NumberFormat = "@"

$ws.Range("K4").Value = "FC-00"
$ws.Range("L4").Value = "999-99"

# Update the active selection to L5 as shown in the diff
$ws.Range("L5").Select()

$wb.Save()
